$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("D2").Value = 0.0294
$ws.Range("D3").Value = 0.0294

$ws.Range("G2").Value = -0.01582278481012658
$ws.Range("G3").Value = -0.01582278481012658

$ws.Range("H2").Value = -0.01582278481012658
$ws.Range("H3").Value = -0.01582278481012658

$ws.Range("I2").Value = -0.2689873417721519
$ws.Range("I3").Value = -0.2689873417721519

$ws.Range("J2").Value = -0.2689873417721519
$ws.Range("J3").Value = -0.2689873417721519

$ws.Range("K2").Value = -0.107
$ws.Range("K3").Value = -0.107

$ws.Range("L2").Value = -0.3386075949367088
$ws.Range("L3").Value = -0.3386075949367088

$ws.Range("U2").Value = 0
$ws.Range("U3").Value = 0

$ws.Range("V2").Value = 0
$ws.Range("V3").Value = 0

$ws.Range("W2").Value = -0.02061657032755299
$ws.Range("W3").Value = -0.02061657032755299

$ws.Range("X2").Value = 0.08291513563551006
$ws.Range("X3").Value = 0.08291513563551006

$ws.Range("Y2").Value = -0.103531705963063
$ws.Range("Y3").Value = -0.103531705963063

$ws.Range("Z2").Value = 0.06089805357486992
$ws.Range("Z3").Value = 0.06089805357486992

$ws.Range("AA2").Value = -0.01638080555020235
$ws.Range("AA3").Value = -0.01638080555020235

$ws.Range("AB2").Value = 0.08291513563551006
$ws.Range("AB3").Value = 0.08291513563551006

$ws.Range("AC2").Value = -0.09929594118571242
$ws.Range("AC3").Value = -0.09929594118571242

$ws.Range("AG2").Value = 0
$ws.Range("AG3").Value = 0

$ws.Range("AJ2").Value = 0
$ws.Range("AJ3").Value = 0

$ws.Range("AK2").Value = 0
$ws.Range("AK3").Value = 0

$ws.Range("AP2").Value = -0
$ws.Range("AP3").Value = -0
